# Vaccine workbook cleanup:
#  - Strip the footnote markers like " [1]" / " [5, 6]" that were appended
#    to vaccine names (the surrounding text/spacing is left untouched).
#  - Collapse the manual line-breaks inside cell text ("`n") into a single
#    space so multi-line labels read as one line.
#  - Fix one data-entry mistake on the "Adult Vaccine " sheet: row 22's NDC
#    number had been typed into the BrandName column (B22) instead of the
#    NDC column (C22), leaving the real NDC cell blank.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $val = $cell.Value2

            if ($val -is [string]) {
                $newVal = $val -replace '\[\d+(,\s*\d+)*\]', ''
                $newVal = $newVal -replace "`n", ' '

                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}

# One-off manual fix mentioned in the commit: the NDC for the plain
# "Tetanus and Diphtheria Toxoids" row had ended up in the BrandName cell.
$adultWs = $wb.Worksheets.Item("Adult Vaccine ")
$adultWs.Cells.Item(22, 3).Value = "00006-4133-41"
$adultWs.Cells.Item(22, 2).ClearContents()
